# Weekly update: prepend a new pair of rows (Primera/Segunda) for the
# latest reporting date, pushing all existing data rows down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 621:622 (shifts rows 621.. down by 2,
# carrying the row-621 formatting, including the date number format).
$ws.Range("A621:R622").Insert()

# New row 621 - "Primera" quality for the new week (2023-07-24 / serial 45131)
$ws.Range("A621").Value = 8
$ws.Range("B621").Value = "Terminal La Palmera de La Serena"
$ws.Range("C621").Value = "Coquimbo"
$ws.Range("D621").Value = 45131
$ws.Range("E621").Value = 4
$ws.Range("F621").Value = 100112009
$ws.Range("G621").Value = "Acelga"
$ws.Range("H621").Value = "Sin especificar"
$ws.Range("I621").Value = "Primera"
$ws.Range("J621").Value = 1600
$ws.Range("K621").Value = 550
$ws.Range("L621").Value = 600
$ws.Range("M621").Value = 575
$ws.Range("N621").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O621").Value = "Provincia del Elquí"
$ws.Range("P621").Value = 288
$ws.Range("Q621").Value = 2
$ws.Range("R621").Value = "Hortaliza"

# New row 622 - "Segunda" quality for the new week (2023-07-24 / serial 45131)
$ws.Range("A622").Value = 8
$ws.Range("B622").Value = "Terminal La Palmera de La Serena"
$ws.Range("C622").Value = "Coquimbo"
$ws.Range("D622").Value = 45131
$ws.Range("E622").Value = 4
$ws.Range("F622").Value = 100112009
$ws.Range("G622").Value = "Acelga"
$ws.Range("H622").Value = "Sin especificar"
$ws.Range("I622").Value = "Segunda"
$ws.Range("J622").Value = 800
$ws.Range("K622").Value = 450
$ws.Range("L622").Value = 500
$ws.Range("M622").Value = 475
$ws.Range("N622").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O622").Value = "Provincia del Elquí"
$ws.Range("P622").Value = 238
$ws.Range("Q622").Value = 2
$ws.Range("R622").Value = "Hortaliza"
